$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 24, shifting rows 24-27 down to 25-28
$ws.Rows(24).Insert()

# Fill in the new row 24 with data
$ws.Range("A24").Value = 3
$ws.Range("B24").Value = "Femacal de La Calera"
$ws.Range("C24").Value = "Coquimbo"
$ws.Range("D24").Value = 45006
$ws.Range("D24").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E24").Value = 5
$ws.Range("F24").Value = "Fruta"
$ws.Range("G24").Value = 100104
$ws.Range("H24").Value = "Frutos de pepita"
$ws.Range("I24").Value = 100104001
$ws.Range("J24").Value = "Granada"
$ws.Range("K24").Value = "Wonderfull"
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 40
$ws.Range("N24").Value = 16000
$ws.Range("O24").Value = 16000
$ws.Range("P24").Value = 16000
$ws.Range("Q24").Value = "$/caja 14 kilos empedrada"
$ws.Range("R24").Value = "Provincia del Elquí"
$ws.Range("S24").Value = 1143
$ws.Range("T24").Value = 14
